# Update countries & provincias Spain
# Applies the refreshed COVID-19 case counts and re-orders the three
# country rows whose ranking (by total cases, column B) changed, while
# keeping the still-correct row data attached to the countries that own it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple data refreshes (country/name unchanged, only counts updated) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 5823877
$ws.Range("C4").Value = 27150
$ws.Range("D4").Value = 3134584
$ws.Range("E4").Value = 2509556
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 537
$ws.Range("H4").Value = 179737

# Row 6: India
$ws.Range("B6").Value = 3043436
$ws.Range("C6").Value = 70068
$ws.Range("D6").Value = 2279900
$ws.Range("E6").Value = 706690
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 918
$ws.Range("H6").Value = 56846

# Row 22: Francia
$ws.Range("B22").Value = 238002
$ws.Range("C22").Value = 3602
$ws.Range("D22").Value = 84950
$ws.Range("E22").Value = 122540
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = 30512

# Row 127: Mali
$ws.Range("B127").Value = 2699
$ws.Range("C127").Value = 11
$ws.Range("D127").Value = 2010
$ws.Range("E127").Value = 564
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 125

# Row 161: Republica del Chad
$ws.Range("B161").Value = 982
$ws.Range("C161").Value = 1
$ws.Range("D161").Value = 869
$ws.Range("E161").Value = 37
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 76

# --- Re-ranked rows ---
# Costa Rica overtakes Moldavia: Costa Rica now gets fresh data in row 65
# and Moldavia (unchanged numbers) drops to row 66.

$ws.Range("A65").Value = "Costa Rica"
$ws.Range("B65").Value = 33084
$ws.Range("C65").Value = 950
$ws.Range("D65").Value = 10372
$ws.Range("E65").Value = 22364
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 8
$ws.Range("H65").Value = 348

$ws.Range("A66").Value = "Moldavia"
$ws.Range("B66").Value = 33072
$ws.Range("C66").Value = 588
$ws.Range("D66").Value = 22683
$ws.Range("E66").Value = 9454
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 6
$ws.Range("H66").Value = 935

# Namibia overtakes Republica de Yibuti (and Malaui): Namibia now gets
# fresh data in row 106; Yibuti and Malaui (unchanged numbers) each drop
# one row.

$ws.Range("A106").Value = "Namibia"
$ws.Range("B106").Value = 5538
$ws.Range("C106").Value = 311
$ws.Range("D106").Value = 2460
$ws.Range("E106").Value = 3032
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 4
$ws.Range("H106").Value = 46

$ws.Range("A107").Value = "Republica de Yibuti"
$ws.Range("B107").Value = 5382
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 5233
$ws.Range("E107").Value = 89
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 60

$ws.Range("A108").Value = "Malaui"
$ws.Range("B108").Value = 5322
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 2929
$ws.Range("E108").Value = 2227
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 166

# Siria overtakes Benin: Siria now gets fresh data in row 136; Benin
# (unchanged numbers) drops to row 137.

$ws.Range("A136").Value = "Siria"
$ws.Range("B136").Value = 2143
$ws.Range("C136").Value = 70
$ws.Range("D136").Value = 490
$ws.Range("E136").Value = 1568
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 2
$ws.Range("H136").Value = 85

$ws.Range("A137").Value = "Benin"
$ws.Range("B137").Value = 2095
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 1705
$ws.Range("E137").Value = 351
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 39

# --- Timestamp header update (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Agosto de 2020 a las 21:41"
